# Apply cryptocurrency price/volume updates per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.669.89'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '3.394.11'
$ws.Range("E3").Value = '  +0.07%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.49'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.91'
$ws.Range("E6").Value = '  -0.32%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.50%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.66'
$ws.Range("E9").Value = '  +0.45%  '
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("E11").Value = '  -2.16%  '
$ws.Range("D12").Value = '3.974.22'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("E13").Value = '  -0.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.20'
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = '3.384.67'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("D17").Value = '61.763.00'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.12'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.04'
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.10'
$ws.Range("E21").Value = '  +2.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.79'
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.549'
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("E24").Value = '  +0.14%  '
$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.197'
$ws.Range("E25").Value = '  +9.48%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000113'
$ws.Range("E26").Value = '  -4.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.997'
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.41'
$ws.Range("E28").Value = '  +1.13%  '
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("E30").Value = '  -0.46%  '
$ws.Range("E31").Value = '  +2.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.32'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.93'
$ws.Range("E34").Value = '  -0.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '169.27'
$ws.Range("E35").Value = '  +1.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.06'
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").Value = '3.428.06'
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0767'
$ws.Range("E39").Value = '  -0.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.79'
$ws.Range("E40").Value = '  -4.92%  '
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.42'
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.66'
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("D45").Value = '2.449.75'
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.76'
$ws.Range("E46").Value = '  -0.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.65'
$ws.Range("E47").Value = '  -2.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0262'
$ws.Range("E49").Value = '  -1.88%  '
$ws.Range("E50").Value = '  -6.16%  '
$ws.Range("E51").Value = '  -1.06%  '
